$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Period:" value cell -> use client timezone aware formatting (joda-style) instead of
# the old String.format based period expression.
$ws.Range("B6").Value = '${from.toString("YYYY.MM.dd HH:mm:ss")+" - "+to.toString("YYYY.MM.dd HH:mm:ss")}'

# Event row "Time" placeholder -> convert event.serverTime into the client timezone using
# joda-time's DateTime constructor before formatting.
$ws.Range("A9").Value = '${new("org.joda.time.DateTime", event.serverTime, timezone).toString("YYYY.MM.dd HH:mm:ss")}'
